$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8204
$ws1.Range("F11").Value = 12
$ws1.Range("F14").Value = 37
$ws1.Range("F16").Value = 570
$ws1.Range("F18").Value = 66
$ws1.Range("F21").Value = 7100
$ws1.Range("F23").Value = 54909
$ws1.Range("F24").Value = 54909
$ws1.Range("F25").Value = 4357
$ws1.Range("F27").Value = 844
$ws1.Range("F28").Value = 409
$ws1.Range("F33").Value = 2894
$ws1.Range("F35").Value = 25
$ws1.Range("F38").Value = 1164
$ws1.Range("F39").Value = 827
$ws1.Range("F40").Value = 146
$ws1.Range("F43").Value = 697
$ws1.Range("F46").Value = 8
$ws1.Range("F47").Value = 143
$ws1.Range("F49").Value = 35
$ws1.Range("F50").Value = 2467

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 46
$ws2.Range("F15").Value = 168
$ws2.Range("F16").Value = 7445
$ws2.Range("F39").Value = 22

# Sheet: 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2274
$ws3.Range("F5").Value = 1536
$ws3.Range("F9").Value = 9321
$ws3.Range("F10").Value = 1637
$ws3.Range("F12").Value = 77
$ws3.Range("F15").Value = 148

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8204
$ws4.Range("F6").Value = 1638
$ws4.Range("F7").Value = 77
$ws4.Range("F13").Value = 46
$ws4.Range("F16").Value = 37
$ws4.Range("F18").Value = 66
$ws4.Range("F20").Value = 54909
$ws4.Range("F23").Value = 844
$ws4.Range("F24").Value = 409
$ws4.Range("F30").Value = 25
$ws4.Range("F32").Value = 1164
$ws4.Range("F35").Value = 146
$ws4.Range("F37").Value = 697
$ws4.Range("F44").Value = 143
$ws4.Range("F46").Value = 35
$ws4.Range("F49").Value = 2467
$ws4.Range("F50").Value = 22

$wb.Save()
